# Fruta / hortaliza, semanal
#
# Insert a new week's worth of price observations (4 rows covering the
# Especial / Primera / Segunda / Tercera quality grades) for
# "Terminal La Palmera de La Serena" - Pina, right before the existing
# row 953, shifting every subsequent row down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 953..956; existing rows 953+ move down to 957+.
$ws.Range("A953:T956").EntireRow.Insert()

# New week's data (columns A..T).
$newData = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44746, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Especial", 216, 21000, 22000, 21500, "`$/caja 10 unidades", "Ecuador", 2150, 10),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44746, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Primera", 216, 21000, 22000, 21500, "`$/caja 12 unidades", "Ecuador", 1792, 12),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44746, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Segunda", 216, 21000, 22000, 21500, "`$/caja 14 unidades", "Ecuador", 1536, 14),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44746, 4, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Tercera", 208, 21000, 22000, 21519, "`$/caja 16 unidades", "Ecuador", 1345, 16)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowNum = 953 + $i
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
